$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired order (language, value) for rows 2..21, sorted by value descending,
# with "Bengali" and "Uzbek" removed from the original data set.
$data = @(
    @("English", 25.11399039441374),
    @("Chinese", 11.62738378438469),
    @("Spanish", 7.486971564569162),
    @("Japanese", 5.819283894129574),
    @("Arabic", 5.095215917791728),
    @("German", 4.928568066176125),
    @("Russian", 3.63058672366909),
    @("Portuguese", 3.473733291049738),
    @("French", 3.102354613874219),
    @("Italian", 2.757308031542069),
    @("Malay-Indonesian", 2.669739947269024),
    @("Korean", 1.653076831763713),
    @("Persian", 1.545525385884489),
    @("Turkish", 1.48050439345805),
    @("Dutch", 1.447744122055965),
    @("Thai", 0.9780226918133414),
    @("Urdu", 0.8626171290788626),
    @("Polish", 0.860320005239489),
    @("Swedish", 0.4842501326030221),
    @("Vietnamese", 0.468776189366164)
)

# Remove the two rows that are dropped from the sheet (Bengali, Uzbek),
# shrinking the used range from A1:B23 down to A1:B21.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

# Write the sorted language/value pairs into rows 2..21.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
